$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ('Price') holds text-formatted numbers (e.g. '1.000', '22.00').
# Temporarily force Text format so Excel doesn't auto-convert the strings to
# numeric values, then restore the original 'General' format afterwards.
$ws.Range("D2:D51").NumberFormat = "@"

# Data rows 2-51: Coin, Link, Price, Volume(1h)
$data = @(
    @('Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '30.727.13', '  +0.61%  '),
    @('Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.888.70', '  +0.23%  '),
    @('TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '0.9996', '  +0.03%  '),
    @('BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '249.63', '  +0.92%  '),
    @('USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.9998', '  +0.10%  '),
    @('XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4757', '  +0.08%  '),
    @('Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2939', '  +1.15%  '),
    @('Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06543', '  +0.17%  '),
    @('Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '22.10', '  +0.17%  '),
    @('TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07735', '  +0.00%  '),
    @('Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.7419', '  -0.43%  '),
    @('Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '96.81', '  -0.82%  '),
    @('WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.887.24', '  +0.14%  '),
    @('Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.244', '  +1.79%  '),
    @('BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '275.72', '  +0.02%  '),
    @('WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '30.796.20', '  +0.97%  '),
    @('Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '13.22', '  -3.05%  '),
    @('ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000007572', '  -0.05%  '),
    @('Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.9996', '  -0.03%  '),
    @('WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.133.98', '  +0.43%  '),
    @('Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '5.340', '  +1.15%  '),
    @('BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '0.9991', '  -0.01%  '),
    @('Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '6.254', '  +0.76%  '),
    @('Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '9.250', '  -1.00%  '),
    @('Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '164.18', '  +0.49%  '),
    @('EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '18.85', '  -0.53%  '),
    @('LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.926', '  -1.20%  '),
    @('Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.345', '  -2.46%  '),
    @('Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.09735', '  -2.43%  '),
    @('PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.507', '  -1.01%  '),
    @('Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.303', '  -0.63%  '),
    @('InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.172', '  +2.30%  '),
    @('Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.04892', '  +1.85%  '),
    @('ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.127', '  -0.29%  '),
    @('ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.6997', '  -0.45%  '),
    @('HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.721', '  -0.01%  '),
    @('VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01911', '  +1.82%  '),
    @('MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.800', '  +2.44%  '),
    @('FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.347', '  +0.15%  '),
    @('Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '75.62', '  +6.05%  '),
    @('RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '2.026', '  +2.99%  '),
    @('TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.4263', '  +0.16%  '),
    @('TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.8418', '  +0.02%  '),
    @('PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '0.9995', '  +0.03%  '),
    @('Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '102.70', '  -0.19%  '),
    @('EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '9.369', '  +0.36%  '),
    @('Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '7.071', '  -0.55%  '),
    @('Elrond', 'https://coinranking.com/coin/omwkOTglq+elrond-egld', '35.69', '  +0.14%  '),
    @('Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '919.20', '  -0.49%  '),
    @('Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.05769', '  +2.04%  ')
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 2).Value = $item[0]
    $ws.Cells.Item($row, 3).Value = $item[1]
    $ws.Cells.Item($row, 4).Value = $item[2]
    $ws.Cells.Item($row, 5).Value = $item[3]
    $row++
}

$ws.Range("D2:D51").NumberFormat = "General"
